# added random generated rec + work on hybrid + small changes overall

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in C22: "Internet of Thingd" -> "Internet of Things"
$c22 = $ws.Range("C22").Value2
$ws.Range("C22").Value = $c22.Replace("Internet of Thingd", "Internet of Things")

# Add new column E: header "Random Recommendation" and set its width
$ws.Range("E1").Value = "Random Recommendation"
# ColumnWidth (COM) differs from the stored OOXML width by the sheet's
# fixed digit-padding offset (~0.83 for this workbook's font), so back it
# out to land on a stored width of exactly 23.
$ws.Columns.Item(5).ColumnWidth = 22.17
